$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G62").Value = 2.02
$ws.Range("H62").Value = 2.03
$ws.Range("I62").Value = 1

$ws.Range("G63").Value = 4.02
$ws.Range("H63").Value = 4.05
$ws.Range("I63").Value = 1

$ws.Range("G64").Value = 6.02
$ws.Range("H64").Value = 6.04
$ws.Range("I64").Value = 1

$ws.Range("F65").Value = 12.5
$ws.Range("G65").Value = 8.03
$ws.Range("H65").Value = 8.01
$ws.Range("I65").Value = 1

$ws.Range("F66").Value = 10
$ws.Range("G66").Value = 10
$ws.Range("H66").Value = 10
$ws.Range("I66").Value = 1

$ws.Range("F69").Select()
